$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: shift day-2 data (rows 98-193) into day-1 slot (rows 2-97).
# The oldest day (25.01.2026) is dropped; 26.01.2026 becomes the new first day.
for ($r = 2; $r -le 97; $r++) {
    $src = $r + 96
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($src, 1).Value2
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($src, 2).Value2
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($src, 3).Value2
}

# Step 2: refresh the "Lookup" labels for the shifted rows (now 26.01.2026).
for ($q = 1; $q -le 96; $q++) {
    $ws.Cells.Item($q + 1, 5).Value = "26.01.2026$q"
}

# Step 3: append the freshly retrained 27.01.2026 day (rows 98-193).
$day3 = @(
    @(46049, 1974.837, 2068),
    @(46049.01041666666, 1954.002, 1511),
    @(46049.02083333334, 1932.784, 1876),
    @(46049.03125, 1925.656, 2019),
    @(46049.04166666666, 1883.287, 1991),
    @(46049.05208333334, 1843.881, 1963),
    @(46049.0625, 1798.173, 1920),
    @(46049.07291666666, 1757.249, 1875),
    @(46049.08333333334, 1703.303, 1892),
    @(46049.09375, 1644.679, 1835),
    @(46049.10416666666, 1592.678, 1735),
    @(46049.11458333334, 1534.322, 1676),
    @(46049.125, 1477.361, 1651),
    @(46049.13541666666, 1457.716, 1647),
    @(46049.14583333334, 1407.104, 1603),
    @(46049.15625, 1358.555, 1539),
    @(46049.16666666666, 1257.145, 1429),
    @(46049.17708333334, 1207.411, 1405),
    @(46049.1875, 1157.217, 1328),
    @(46049.19791666666, 1109.479, 1279),
    @(46049.20833333334, 965.3680000000001, 1220),
    @(46049.21875, 930.788, 1123),
    @(46049.22916666666, 891.154, 1001),
    @(46049.23958333334, 837.203, 977),
    @(46049.25, 743.835, 993),
    @(46049.26041666666, 686.65, 901),
    @(46049.27083333334, 647.857, 798),
    @(46049.28125, 604.904, 716),
    @(46049.29166666666, 568.997, 685),
    @(46049.30208333334, 548.943, 652),
    @(46049.3125, 523.662, 600),
    @(46049.32291666666, 494.949, 520),
    @(46049.33333333334, 444.617, 429),
    @(46049.34375, 425.238, 406),
    @(46049.35416666666, 390.342, 394),
    @(46049.36458333334, 369.207, 356),
    @(46049.375, 316.628, 299),
    @(46049.38541666666, 302.658, 250),
    @(46049.39583333334, 275.073, 219),
    @(46049.40625, 261.815, 223),
    @(46049.41666666666, 238.008, 204),
    @(46049.42708333334, 231.883, 182),
    @(46049.4375, 225.48, 165),
    @(46049.44791666666, 218.87, 0),
    @(46049.45833333334, 206.183, 0),
    @(46049.46875, 201.691, 0),
    @(46049.47916666666, 196.876, 0),
    @(46049.48958333334, 192.095, 0),
    @(46049.5, 187.23, 0),
    @(46049.51041666666, 184.178, 0),
    @(46049.52083333334, 181.47, 0),
    @(46049.53125, 179.028, 0),
    @(46049.54166666666, 175.407, 0),
    @(46049.55208333334, 175.742, 0),
    @(46049.5625, 177.446, 0),
    @(46049.57291666666, 178.5, 0),
    @(46049.58333333334, 181.447, 0),
    @(46049.59375, 181.352, 0),
    @(46049.60416666666, 181.744, 0),
    @(46049.61458333334, 196.884, 0),
    @(46049.625, 188.206, 0),
    @(46049.63541666666, 209.324, 0),
    @(46049.64583333334, 215.563, 0),
    @(46049.65625, 221.56, 0),
    @(46049.66666666666, 216.509, 0),
    @(46049.67708333334, 225.956, 0),
    @(46049.6875, 248.806, 0),
    @(46049.69791666666, 258.029, 0),
    @(46049.70833333334, 260.149, 0),
    @(46049.71875, 277.895, 0),
    @(46049.72916666666, 282.456, 0),
    @(46049.73958333334, 275.149, 0),
    @(46049.75, 295.977, 0),
    @(46049.76041666666, 304.318, 0),
    @(46049.77083333334, 312.819, 0),
    @(46049.78125, 310.149, 0),
    @(46049.79166666666, 332.131, 0),
    @(46049.80208333334, 342.978, 0),
    @(46049.8125, 352.459, 0),
    @(46049.82291666666, 352.315, 0),
    @(46049.83333333334, 377.231, 0),
    @(46049.84375, 373.78, 0),
    @(46049.85416666666, 380.419, 0),
    @(46049.86458333334, 387.249, 0),
    @(46049.875, 396.631, 0),
    @(46049.88541666666, 401.198, 0),
    @(46049.89583333334, 404.451, 0),
    @(46049.90625, 406.218, 0),
    @(46049.91666666666, 407.996, 0),
    @(46049.92708333334, 411.749, 0),
    @(46049.9375, 415.837, 0),
    @(46049.94791666666, 421.176, 0),
    @(46049.95833333334, 0, 0),
    @(46049.96875, 0, 0),
    @(46049.97916666666, 0, 0),
    @(46049.98958333334, 0, 0)
)
for ($i = 0; $i -lt $day3.Length; $i++) {
    $r = $i + 98
    $row = $day3[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $q = $i + 1
    $ws.Cells.Item($r, 5).Value = "27.01.2026$q"
}